# feat: add 2022-Q3 data
#
# Insert a new quarterly sheet "2022-Q3" right after "总计", pushing the
# existing "2022-Q2" / "2022-Q1" / "2021-Q4" sheets one slot later (their
# own contents are untouched - they just shift position). Populate the new
# sheet with the Q3 fund-holding data, and update the "总计" (totals) sheet
# with the new quarter's summary row plus the now-shifted rows.

$wb = $excel.ActiveWorkbook

$totals = $wb.Worksheets.Item(1)
$lastQuarterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- Create the new "2022-Q3" sheet by copying an existing quarter sheet so
# --- it inherits identical formatting/styles, then place it right after
# --- "总计". The sheets after it (2022-Q2, 2022-Q1, 2021-Q4) automatically
# --- shift down by one position and keep their own data untouched.
$lastQuarterSheet.Copy($null, $totals)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Force text (string) formatting on the columns that hold text-typed values
# (fund code, fund size, stock position, position ratio, market value) so
# they don't get silently reinterpreted as numbers (losing e.g. the leading
# zero on fund codes).
$q3.Range("B2:B3").NumberFormat = "@"
$q3.Range("D2:G3").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "012349"
$q3.Range("C2").Value = "天弘恒生科技指数（QDII）C"
$q3.Range("D2").Value = "33.57"
$q3.Range("E2").Value = "92.84"
$q3.Range("F2").Value = "5.91"
$q3.Range("G2").Value = "1.9840"
$q3.Range("H2").Value = 7

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "012348"
$q3.Range("C3").Value = "天弘恒生科技指数（QDII）A"
$q3.Range("D3").Value = "30.64"
$q3.Range("E3").Value = "92.84"
$q3.Range("F3").Value = "5.91"
$q3.Range("G3").Value = "1.8108"
$q3.Range("H3").Value = 7

# Drop back to the default ("Normal") style so these text-forced cells don't
# pick up a stray number-format style index (matches the unstyled cells the
# other quarter sheets use for the same columns).
$q3.Range("B2:B3").Style = "Normal"
$q3.Range("D2:G3").Style = "Normal"

# --- Update the "总计" (totals) sheet: shift the existing rows down one
# --- and add the new 2022-Q3 summary row at the top of the data.
$totals.Range("A4:D4").Copy($totals.Range("A5:D5"))

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 2
$totals.Range("D2").Value = 3.79

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q2"
$totals.Range("C3").Value = 2
$totals.Range("D3").Value = 3.83

$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2022-Q1"
$totals.Range("C4").Value = 2
$totals.Range("D4").Value = 4.2

$totals.Range("A5").Value = 3
$totals.Range("B5").Value = "2021-Q4"
$totals.Range("C5").Value = 2
$totals.Range("D5").Value = 1.2
